# Refresh the crypto price/volume snapshot (GitHub Actions scheduled update).
# Column D ("Price") holds text that often *looks* numeric (e.g. "211.34",
# thousands-dotted "27.514.95"), so we force the Text number format before
# writing and reset the style afterwards — otherwise Excel's COM layer would
# silently coerce parseable values (like "211.34") into real numbers instead
# of leaving them as the original literal text.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "27.514.95"
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = "  -0.90%  "
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "1.620.08"
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = "  -1.71%  "
$ws.Range("E4").Value = "  -0.07%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "211.34"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  -0.97%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "0.526"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  -0.60%  "
$ws.Range("E7").Value = "  -0.07%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "23.17"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "  -1.65%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.263"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  +1.76%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.0612"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  -0.50%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.0887"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  -0.03%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "1.848.66"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  -1.77%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "1.629.76"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  -1.12%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "4.05"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  -0.05%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "0.550"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  -2.03%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "65.08"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  +0.75%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "27.475.54"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "  -0.98%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "233.21"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  +0.17%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "0.0₃0719"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  -0.92%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "7.53"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  -1.20%  "
$ws.Range("E21").Value = "  -0.10%  "
$ws.Range("B22").Value = "Uniswap"
$ws.Range("C22").Value = "https://coinranking.com/coin/_H5FVG9iW+uniswap-uni"
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "4.32"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  +0.40%  "
$ws.Range("B23").Value = "Avalanche"
$ws.Range("C23").Value = "https://coinranking.com/coin/dvUj0CzDZ+avalanche-avax"
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "10.20"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  +0.35%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "2.07"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  +6.29%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "149.91"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  +0.09%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "6.86"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  -1.14%  "
$ws.Range("E27").Value = "  -0.32%  "
$ws.Range("E28").Value = "  -0.12%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "15.54"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  -1.12%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "1.17"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  -1.38%  "
$ws.Range("E31").Value = "  -1.00%  "
$ws.Range("E32").Value = "  -0.81%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "1.481.95"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  +3.19%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "3.08"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  -2.50%  "
$ws.Range("E35").Value = "  -2.98%  "
$ws.Range("E36").Value = "  -0.69%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "0.959"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  +8.18%  "
$ws.Range("E38").Value = "  -0.30%  "
$ws.Range("E39").Value = "  -2.92%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.867"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  -1.99%  "
$ws.Range("E41").Value = "  -0.09%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "1.02"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  -2.07%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "67.62"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  +1.06%  "
$ws.Range("B44").Value = "MXToken"
$ws.Range("C44").Value = "https://coinranking.com/coin/QUC5kVAxSoB-+mxtoken-mx"
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "2.20"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  -2.37%  "
$ws.Range("B45").Value = "RocketPoolETH"
$ws.Range("C45").Value = "https://coinranking.com/coin/QJZRUGyNI+rocketpooleth-reth"
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "1.759.10"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  -1.80%  "
$ws.Range("B46").Value = "FraxShare"
$ws.Range("C46").Value = "https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs"
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "5.22"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  -6.53%  "
$ws.Range("B47").Value = "RenderToken"
$ws.Range("C47").Value = "https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr"
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "1.74"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  -0.73%  "
$ws.Range("B48").Value = "Quant"
$ws.Range("C48").Value = "https://coinranking.com/coin/bauj_21eYVwso+quant-qnt"
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "87.00"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  +1.93%  "
$ws.Range("B49").Value = "Algorand"
$ws.Range("C49").Value = "https://coinranking.com/coin/TpHE2IShQw-sJ+algorand-algo"
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "0.101"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  +1.14%  "
$ws.Range("B50").Value = "Cronos"
$ws.Range("C50").Value = "https://coinranking.com/coin/65PHZTpmE55b+cronos-cro"
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "0.0502"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  -0.18%  "
$ws.Range("B51").Value = "EnergySwap"
$ws.Range("C51").Value = "https://coinranking.com/coin/SbWqqTui-+energyswap-ens"
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "7.64"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  -1.66%  "
